# The document has two logo pictures that repeat across the primary and
# first-page headers/footers:
#   - headers: the BTEC logo (currently displayed/internal picture name
#     "image1.jpg")
#   - footers: the Pearson Edexcel logo (currently displayed/internal
#     picture name "image2.png")
#
# The authored change simply swaps each picture's display "Name" (the
# value you see/edit in the Selection Pane, backed by <wp:docPr name="...">
# in the OOXML) to the other sibling's label:
#   headers: image1.jpg -> image2.jpg
#   footers: image2.png -> image1.png
#
# Inline pictures (wdInlineShapePicture) don't expose a settable .Name on
# the Word InlineShape object itself, so we briefly promote each inline
# picture to a floating Shape (ConvertToShape), rename it there (where
# .Name is writable), then convert it straight back to an inline picture
# so layout/wrapping is unaffected.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $floatingShape = $inlineShape.ConvertToShape()
    $floatingShape.Name = $newName
    $floatingShape.ConvertToInlineShape() | Out-Null
}

# Headers: BTEC logo, image1.jpg -> image2.jpg
Rename-InlinePicture $sec.Headers.Item(1).Range "image2.jpg"
Rename-InlinePicture $sec.Headers.Item(2).Range "image2.jpg"

# Footers: Pearson Edexcel logo, image2.png -> image1.png
Rename-InlinePicture $sec.Footers.Item(1).Range "image1.png"
Rename-InlinePicture $sec.Footers.Item(2).Range "image1.png"
